$d = $word.ActiveDocument

# 1) Fix "contactor" -> "contractor"
$d.Content.Find.Execute("contactor", $false, $false, $false, $false, $false,
                         $true, 1, $false, "contractor", 2)

# 2) Split the run containing the stray trailing "w" and remove it.
#    "Worked for e-learning, e-commerce and banking domains.w"
#    becomes two runs: "Worked for e-learning, e-commerce and banking "
#    and "domains." (trailing stray "w" removed).
$d.Content.Find.Execute("Worked for e-learning, e-commerce and banking domains.w", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Worked for e-learning, e-commerce and banking domains.", 2)
